# Update the public EPEX Spot prices workbook:
#  - "Prix Spot": insert a new "06-dec" column (before the "01-oct." block)
#    with no data yet ("-" placeholders) for every hourly row.
#  - "Gaz" / "CO2": append the next day's closing price as a new row.

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert column EI (06-dec), shifting 01-oct..31-oct right ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

$wsSpot.Columns("EI:EI").Insert()
$wsSpot.Range("EI1").Value = "06-dec"
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 139).Value = "-"
}

# --- Sheet "Gaz": append 2025-12-04 price ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A169").NumberFormat = "@"
$wsGaz.Range("A169").Value = "2025-12-04"
$wsGaz.Range("A169").Style = "Normal"
$wsGaz.Range("B169").Value = 25.95

# --- Sheet "CO2": append 2025-12-04 price ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A169").NumberFormat = "@"
$wsCo2.Range("A169").Value = "2025-12-04"
$wsCo2.Range("A169").Style = "Normal"
$wsCo2.Range("B169").Value = 82.5
